$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: michael@example.com duplicate -> jane@example.com / xyz987hashed
$ws.Range("A4").Value = "jane@example.com"
$ws.Range("B4").Value = "xyz987hashed"

# Row 5: invalidUser row - Actual Result Pass -> Fail, Status Fail -> Pass
$ws.Range("D5").Value = "Fail"
$ws.Range("E5").Value = "Pass"

# Row 6: validUser -> john@example.com, invalidPass stays, Actual Result Pass -> Fail, Status Fail -> Pass
$ws.Range("A6").Value = "john@example.com"
$ws.Range("B6").Value = "invalidPass"
$ws.Range("D6").Value = "Fail"
$ws.Range("E6").Value = "Pass"

# Remove old rows 7 and 8 (shrink the used range to A1:E6)
$ws.Rows("7:8").Delete()
